$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.997.77'
$ws.Range("E2").Value = '  -1.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.824.60'
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.26'
$ws.Range("E5").Value = '  -1.22%  '

$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4254'
$ws.Range("E7").Value = '  -1.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3657'
$ws.Range("E8").Value = '  -1.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07237'
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8432'
$ws.Range("E10").Value = '  -2.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.58'
$ws.Range("E11").Value = '  -2.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.832.91'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.653'
$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07048'
$ws.Range("E14").Value = '  -0.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.287'
$ws.Range("E15").Value = '  -1.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.62'
$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008757'
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.88'
$ws.Range("E20").Value = '  -2.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.058.11'
$ws.Range("E21").Value = '  -1.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.132'
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.86'
$ws.Range("E23").Value = '  -0.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.050.03'
$ws.Range("E24").Value = '  -0.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.981'
$ws.Range("E25").Value = '  -0.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.11'
$ws.Range("E26").Value = '  -1.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.216'
$ws.Range("E27").Value = '  +3.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.19'
$ws.Range("E28").Value = '  -1.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.229'
$ws.Range("E29").Value = '  -1.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.81'
$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08706'
$ws.Range("E31").Value = '  -1.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.177'
$ws.Range("E32").Value = '  -2.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7389'
$ws.Range("E33").Value = '  -3.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.900'
$ws.Range("E34").Value = '  -0.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.419'
$ws.Range("E35").Value = '  -1.42%  '

$ws.Range("E36").Value = '  -0.32%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.091'
$ws.Range("E37").Value = '  -2.70%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01940'
$ws.Range("E38").Value = '  -1.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05216'
$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.251'
$ws.Range("E40").Value = '  +1.08%  '

$ws.Range("E41").Value = '  -0.30%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5130'
$ws.Range("E42").Value = '  +0.93%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1689'
$ws.Range("E43").Value = '  +0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.550'
$ws.Range("E44").Value = '  -1.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.52'
$ws.Range("E45").Value = '  -0.34%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4734'
$ws.Range("E46").Value = '  -0.14%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.946'
$ws.Range("E47").Value = '  +6.24%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.74'
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9996'
$ws.Range("E49").Value = '  -0.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06328'
$ws.Range("E50").Value = '  -1.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.655'
$ws.Range("E51").Value = '  -0.70%  '
